$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.704.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.67%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.900.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.26%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.91%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.86%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.504"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.69%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.899.48"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.22%  "

$ws.Range("E11").Value = "  -4.56%  "

$ws.Range("E12").Value = "  -4.38%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000226"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.29%  "

$ws.Range("E15").Value = "  +1.84%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.382.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.23%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.644.56"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.72%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.900.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "428.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.97%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.682"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.47%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("E29").Value = "  -1.78%  "

$ws.Range("E30").Value = "  +0.01%  "

$ws.Range("E31").Value = "  -3.44%  "

$ws.Range("E32").Value = "  -7.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.48"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.35%  "

$ws.Range("E34").Value = "  -3.46%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0849"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.25%  "

$ws.Range("E36").Value = "  -3.57%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.08%  "

$ws.Range("E38").Value = "  -5.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.23%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.125"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.99"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.89%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.61"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.02%  "

$ws.Range("E43").Value = "  -5.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.56%  "

$ws.Range("E45").Value = "  -3.26%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "371.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.696.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.70%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.95%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.54%  "

$ws.Range("E51").Value = "  -2.55%  "

